$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57, pushing the existing rows 57-60 down to 58-61
$ws.Rows(57).Insert()

# Fill in the newly-inserted row 57 with the new weekly record
$ws.Cells.Item(57, 1).Value2 = 2
$ws.Cells.Item(57, 2).Value2 = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(57, 3).Value2 = "Coquimbo"
$ws.Cells.Item(57, 4).Value2 = 44769
$ws.Cells.Item(57, 5).Value2 = 4
$ws.Cells.Item(57, 6).Value2 = 100112026
$ws.Cells.Item(57, 7).Value2 = "Haba"
$ws.Cells.Item(57, 8).Value2 = "Sin especificar"
$ws.Cells.Item(57, 9).Value2 = "Primera"
$ws.Cells.Item(57, 10).Value2 = 1300
$ws.Cells.Item(57, 11).Value2 = 7000
$ws.Cells.Item(57, 12).Value2 = 8000
$ws.Cells.Item(57, 13).Value2 = 7500
$ws.Cells.Item(57, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(57, 16).Value2 = 300
$ws.Cells.Item(57, 17).Value2 = 25
$ws.Cells.Item(57, 18).Value2 = "Hortaliza"
